# Inserts a new price-report row at row 366 (Florida King, Provincia de Limarí,
# $/bandeja 10 kilos granel) and shifts all subsequent rows down by one,
# growing the used range from A1:T438 to A1:T439.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("366:366").Insert()

$ws.Cells.Item(366, 1).Value = 10
$ws.Cells.Item(366, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(366, 3).Value = "La Araucanía"
$ws.Cells.Item(366, 4).Value = 45218
$ws.Cells.Item(366, 5).Value = 9
$ws.Cells.Item(366, 6).Value = "Fruta"
$ws.Cells.Item(366, 7).Value = 100103
$ws.Cells.Item(366, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(366, 9).Value = 100103004
$ws.Cells.Item(366, 10).Value = "Durazno"
$ws.Cells.Item(366, 11).Value = "Florida King"
$ws.Cells.Item(366, 12).Value = "Primera"
$ws.Cells.Item(366, 13).Value = 250
$ws.Cells.Item(366, 14).Value = 24000
$ws.Cells.Item(366, 15).Value = 24000
$ws.Cells.Item(366, 16).Value = 24000
$ws.Cells.Item(366, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(366, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(366, 19).Value = 2400
$ws.Cells.Item(366, 20).Value = 10
